$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFIDF1")

# Clear the "one seed only" placeholder note in B1 and J1 (style is retained)
$ws.Range("B1").Value = ""
$ws.Range("J1").Value = ""

# Replace the single-seed numeric results with aggregated mean +/- std text results
# for the new TFIDF (1,1) run
$ws.Range("B3").Value = "0.7443±0.0027"
$ws.Range("C3").Value = "0.7035±0.0056"
$ws.Range("D3").Value = "0.5984±0.0046"
$ws.Range("E3").Value = "0.4775±0.0062"
$ws.Range("F3").Value = "0.8217±0.0017"
$ws.Range("G3").Value = "0.6467±0.0044"
$ws.Range("I3").Value = "0.7494±0.0017"
$ws.Range("J3").Value = "0.7139±0.0031"
$ws.Range("K3").Value = "0.6035±0.003"
$ws.Range("L3").Value = "0.4882±0.0041"
$ws.Range("M3").Value = "0.8285±0.0018"
$ws.Range("N3").Value = "0.6541±0.0025"
$ws.Range("B4").Value = "0.7355±0.0027"
$ws.Range("C4").Value = "0.6846±0.006"
$ws.Range("D4").Value = "0.5875±0.0062"
$ws.Range("E4").Value = "0.4589±0.0061"
$ws.Range("F4").Value = "0.8055±0.0033"
$ws.Range("G4").Value = "0.6323±0.0036"
$ws.Range("I4").Value = "0.7348±0.0025"
$ws.Range("J4").Value = "0.6757±0.0038"
$ws.Range("K4").Value = "0.5972±0.004"
$ws.Range("L4").Value = "0.4623±0.0056"
$ws.Range("M4").Value = "0.7972±0.0016"
$ws.Range("N4").Value = "0.634±0.0037"
$ws.Range("B5").Value = "0.7158±0.0034"
$ws.Range("C5").Value = "0.6774±0.0083"
$ws.Range("D5").Value = "0.5302±0.005"
$ws.Range("E5").Value = "0.401±0.0083"
$ws.Range("F5").Value = "0.8058±0.0038"
$ws.Range("G5").Value = "0.5948±0.0059"
$ws.Range("I5").Value = "0.7194±0.003"
$ws.Range("J5").Value = "0.6787±0.0071"
$ws.Range("K5").Value = "0.5415±0.0033"
$ws.Range("L5").Value = "0.412±0.0067"
$ws.Range("M5").Value = "0.8102±0.0041"
$ws.Range("N5").Value = "0.6024±0.0042"
$ws.Range("B6").Value = "0.7237±0.0028"
$ws.Range("C6").Value = "0.6884±0.0075"
$ws.Range("D6").Value = "0.5469±0.0037"
$ws.Range("E6").Value = "0.422±0.0065"
$ws.Range("F6").Value = "0.8115±0.0047"
$ws.Range("G6").Value = "0.6096±0.0049"
$ws.Range("I6").Value = "0.7265±0.002"
$ws.Range("J6").Value = "0.6937±0.0047"
$ws.Range("K6").Value = "0.5513±0.0045"
$ws.Range("L6").Value = "0.4288±0.0052"
$ws.Range("M6").Value = "0.8183±0.0021"
$ws.Range("N6").Value = "0.6143±0.0036"
$ws.Range("B7").Value = "0.7087±0.0026"
$ws.Range("C7").Value = "0.7038±0.0075"
$ws.Range("D7").Value = "0.4994±0.0042"
$ws.Range("E7").Value = "0.3777±0.007"
$ws.Range("F7").Value = "0.8036±0.0027"
$ws.Range("G7").Value = "0.5843±0.005"
$ws.Range("I7").Value = "0.7117±0.0021"
$ws.Range("J7").Value = "0.7152±0.0073"
$ws.Range("K7").Value = "0.5018±0.0038"
$ws.Range("L7").Value = "0.3859±0.0058"
$ws.Range("M7").Value = "0.8068±0.0033"
$ws.Range("N7").Value = "0.5897±0.0045"
$ws.Range("B8").Value = "0.7386±0.0021"
$ws.Range("C8").Value = "0.6986±0.003"
$ws.Range("D8").Value = "0.5856±0.0037"
$ws.Range("E8").Value = "0.4627±0.0052"
$ws.Range("F8").Value = "0.8342±0.0025"
$ws.Range("G8").Value = "0.6371±0.0032"
$ws.Range("I8").Value = "0.7404±0.002"
$ws.Range("J8").Value = "0.7061±0.0048"
$ws.Range("K8").Value = "0.585±0.0036"
$ws.Range("L8").Value = "0.4658±0.0046"
$ws.Range("M8").Value = "0.8397±0.0029"
$ws.Range("N8").Value = "0.6398±0.0036"
$ws.Range("B9").Value = "0.7226±0.0104"
$ws.Range("C9").Value = "0.6538±0.0203"
$ws.Range("D9").Value = "0.6142±0.0089"
$ws.Range("E9").Value = "0.4546±0.0138"
$ws.Range("F9").Value = "0.8198±0.0058"
$ws.Range("G9").Value = "0.6332±0.0083"
$ws.Range("I9").Value = "0.7297±0.0083"
$ws.Range("J9").Value = "0.6666±0.0196"
$ws.Range("K9").Value = "0.6169±0.0061"
$ws.Range("L9").Value = "0.4645±0.0109"
$ws.Range("M9").Value = "0.8265±0.0047"
$ws.Range("N9").Value = "0.6406±0.0071"
$ws.Range("B10").Value = "0.722±0.0056"
$ws.Range("C10").Value = "0.6884±0.0068"
$ws.Range("D10").Value = "0.5462±0.017"
$ws.Range("E10").Value = "0.4183±0.0169"
$ws.Range("F10").Value = "0.8137±0.0047"
$ws.Range("G10").Value = "0.6089±0.0097"
$ws.Range("I10").Value = "0.7251±0.0053"
$ws.Range("J10").Value = "0.6939±0.0064"
$ws.Range("K10").Value = "0.5505±0.0182"
$ws.Range("L10").Value = "0.4259±0.0165"
$ws.Range("M10").Value = "0.8153±0.0063"
$ws.Range("N10").Value = "0.6138±0.0098"
$ws.Range("B11").Value = "0.6467±0.0228"
$ws.Range("C11").Value = "0.5379±0.0461"
$ws.Range("D11").Value = "0.4782±0.0158"
$ws.Range("E11").Value = "0.2656±0.0433"
$ws.Range("F11").Value = "0.6254±0.0641"
$ws.Range("G11").Value = "0.506±0.0284"
$ws.Range("I11").Value = "0.6516±0.0243"
$ws.Range("J11").Value = "0.5467±0.0481"
$ws.Range("K11").Value = "0.489±0.0183"
$ws.Range("L11").Value = "0.2806±0.0471"
$ws.Range("M11").Value = "0.6328±0.064"
$ws.Range("N11").Value = "0.516±0.0308"
$ws.Range("B12").Value = "0.6436±0.027"
$ws.Range("C12").Value = "0.544±0.0476"
$ws.Range("D12").Value = "0.5201±0.0082"
$ws.Range("E12").Value = "0.307±0.0377"
$ws.Range("F12").Value = "0.658±0.0525"
$ws.Range("G12").Value = "0.5312±0.0253"
$ws.Range("I12").Value = "0.6498±0.0267"
$ws.Range("J12").Value = "0.5526±0.0498"
$ws.Range("K12").Value = "0.5262±0.0097"
$ws.Range("L12").Value = "0.3171±0.039"
$ws.Range("M12").Value = "0.6622±0.0539"
$ws.Range("N12").Value = "0.5385±0.0271"
$ws.Range("B13").Value = "0.7322±0.0083"
$ws.Range("C13").Value = "0.6833±0.0091"
$ws.Range("D13").Value = "0.5826±0.029"
$ws.Range("E13").Value = "0.4512±0.0256"
$ws.Range("F13").Value = "0.8286±0.0081"
$ws.Range("G13").Value = "0.6284±0.0151"
$ws.Range("I13").Value = "0.7345±0.0089"
$ws.Range("J13").Value = "0.6874±0.0082"
$ws.Range("K13").Value = "0.5837±0.0295"
$ws.Range("L13").Value = "0.4556±0.0268"
$ws.Range("M13").Value = "0.8334±0.0083"
$ws.Range("N13").Value = "0.6308±0.0157"
$ws.Range("B14").Value = "0.7277±0.005"
$ws.Range("C14").Value = "0.6824±0.0077"
$ws.Range("D14").Value = "0.5614±0.0104"
$ws.Range("E14").Value = "0.435±0.0126"
$ws.Range("F14").Value = "0.8217±0.0056"
$ws.Range("G14").Value = "0.616±0.0093"
$ws.Range("I14").Value = "0.729±0.005"
$ws.Range("J14").Value = "0.6877±0.0085"
$ws.Range("K14").Value = "0.5602±0.0103"
$ws.Range("L14").Value = "0.4371±0.0124"
$ws.Range("M14").Value = "0.824±0.0071"
$ws.Range("N14").Value = "0.6174±0.0095"

# Make TFIDF1 the active sheet/tab and set its selection, mirroring BOW1's
# selection losing "tabSelected" in favor of this sheet.
$ws.Activate()
$ws.Range("Q17").Select()
